$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 3.65
$ws.Range("AI2").Value = 40

# Row 3
$ws.Range("AB3").Value = 980
$ws.Range("AC3").Value = 970

# Row 4
$ws.Range("I4").Value = 4.9
$ws.Range("K4").Value = 4.2
$ws.Range("S4").Value = 2.96
$ws.Range("T4").Value = 1.72

# Row 5
$ws.Range("L5").Value = 1.47
$ws.Range("M5").Value = 1.12
$ws.Range("Y5").Value = 19
$ws.Range("AC5").Value = 10
$ws.Range("AF5").Value = 9.4
$ws.Range("AJ5").Value = 20
$ws.Range("AK5").Value = 28

# Row 6
$ws.Range("X6").Value = 8
$ws.Range("AJ6").Value = 48
$ws.Range("AN6").Value = 55

# Row 8
$ws.Range("G8").Value = 2.08
$ws.Range("H8").Value = 4.3
$ws.Range("W8").Value = 1.93

# Row 9
$ws.Range("F9").Value = 1.93
$ws.Range("G9").Value = 2.14
$ws.Range("I9").Value = 6.4
$ws.Range("J9").Value = 2.76
$ws.Range("K9").Value = 3.7
$ws.Range("L9").Value = 1.65
$ws.Range("V9").Value = 1.2
$ws.Range("W9").Value = 1.88

# Row 10
$ws.Range("J10").Value = 3.95

# Row 11
$ws.Range("F11").Value = 1.81
$ws.Range("G11").Value = 2.04
$ws.Range("H11").Value = 4.2
$ws.Range("I11").Value = 5.2
$ws.Range("V11").Value = 1.24
$ws.Range("W11").Value = 1.96
$ws.Range("X11").Value = 19.5
$ws.Range("Y11").Value = 21
$ws.Range("AD11").Value = 980
$ws.Range("AE11").Value = 70
$ws.Range("AG11").Value = 12.5
$ws.Range("AJ11").Value = 980

# Row 12
$ws.Range("G12").Value = 10.5
$ws.Range("I12").Value = 1.5
$ws.Range("L12").Value = 1.28
$ws.Range("N12").Value = 4.9
$ws.Range("O12").Value = 1.2
$ws.Range("P12").Value = 2.36
$ws.Range("Q12").Value = 1.6
$ws.Range("R12").Value = 1.55
$ws.Range("S12").Value = 2.38
$ws.Range("T12").Value = 1.83
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 2.98
$ws.Range("W12").Value = 1.12
$ws.Range("X12").Value = 29
$ws.Range("Y12").Value = 12
$ws.Range("AA12").Value = 13.5
$ws.Range("AG12").Value = 32
$ws.Range("AH12").Value = 24
$ws.Range("AI12").Value = 34

# Row 19
$ws.Range("X19").Value = 18.5
$ws.Range("AI19").Value = 28
$ws.Range("AL19").Value = 40

# Row 21
$ws.Range("F21").Value = 2.06
$ws.Range("G21").Value = 2.56
$ws.Range("H21").Value = 3.45
$ws.Range("I21").Value = 4.7
$ws.Range("J21").Value = 2.84
$ws.Range("K21").Value = 4.5
$ws.Range("O21").Value = 1.48
$ws.Range("Q21").Value = 1.01
$ws.Range("V21").Value = 1.27
$ws.Range("W21").Value = 1.64

# Row 22
$ws.Range("I22").Value = 3.35
$ws.Range("V22").Value = 1.42
$ws.Range("AE22").Value = 48
$ws.Range("AI22").Value = 65
$ws.Range("AM22").Value = 170
$ws.Range("AN22").Value = 42

# Row 23
$ws.Range("G23").Value = 4.3
$ws.Range("Q23").Value = 1.8
$ws.Range("W23").Value = 1.31
$ws.Range("AB23").Value = 1000
$ws.Range("AE23").Value = 30
$ws.Range("AF23").Value = 1000
$ws.Range("AI23").Value = 44

# Row 24
$ws.Range("Q24").Value = 3.25
$ws.Range("AD24").Value = 14
